$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("其他控制要求")
$ws3.Range("E3:F3").Value = "boolean"
$ws3.Range("E3:F3").WrapText = $true
$ws3.Range("E3:F3").HorizontalAlignment = -4108
$ws3.Range("E3:F3").WrapText = $false
